$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.090.83"
$ws.Range("E2").Value = "  -2.39%  "

$ws.Range("D3").Value = "1.642.95"
$ws.Range("E3").Value = "  -2.09%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'308.74"
$ws.Range("E5").Value = "  -1.57%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").Value = "'0.3944"
$ws.Range("E7").Value = "  +0.71%  "

$ws.Range("D8").Value = "'0.3863"
$ws.Range("E8").Value = "  -2.34%  "

$ws.Range("D9").Value = "'1.004"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Value = "'1.365"
$ws.Range("E10").Value = "  -2.56%  "

$ws.Range("D11").Value = "'49.67"
$ws.Range("E11").Value = "  -4.34%  "

$ws.Range("D12").Value = "'0.08564"
$ws.Range("E12").Value = "  -0.70%  "

$ws.Range("D13").Value = "'23.56"
$ws.Range("E13").Value = "  -6.10%  "

$ws.Range("D14").Value = "'7.084"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("D15").Value = "'0.00001284"
$ws.Range("E15").Value = "  -2.35%  "

$ws.Range("D16").Value = "'7.508"
$ws.Range("E16").Value = "  -3.10%  "

$ws.Range("D17").Value = "1.662.29"
$ws.Range("E17").Value = "  -2.42%  "

$ws.Range("D18").Value = "'93.88"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").Value = "'0.06911"
$ws.Range("E19").Value = "  -2.61%  "

$ws.Range("D20").Value = "'20.30"
$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("D21").Value = "'6.922"
$ws.Range("E21").Value = "  -1.56%  "

$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").Value = "'13.59"
$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("D24").Value = "24.105.87"
$ws.Range("E24").Value = "  -2.35%  "

$ws.Range("D25").Value = "'2.420"
$ws.Range("E25").Value = "  +2.73%  "

$ws.Range("D26").Value = "'2.872"
$ws.Range("E26").Value = "  +4.10%  "

$ws.Range("D27").Value = "'22.18"
$ws.Range("E27").Value = "  -5.52%  "

$ws.Range("D28").Value = "'157.94"
$ws.Range("E28").Value = "  -2.96%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'140.28"
$ws.Range("E29").Value = "  -5.70%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'8.150"
$ws.Range("E30").Value = "  +3.64%  "

$ws.Range("D31").Value = "'5.251"
$ws.Range("E31").Value = "  -9.59%  "

$ws.Range("D32").Value = "'2.521"
$ws.Range("E32").Value = "  +6.10%  "

$ws.Range("D33").Value = "1.832.85"
$ws.Range("E33").Value = "  -2.14%  "

$ws.Range("D34").Value = "'0.08094"
$ws.Range("E34").Value = "  -3.11%  "

$ws.Range("D35").Value = "'6.761"
$ws.Range("E35").Value = "  -2.34%  "

$ws.Range("D36").Value = "'0.02912"
$ws.Range("E36").Value = "  -5.00%  "

$ws.Range("D37").Value = "'0.9656"
$ws.Range("E37").Value = "  -3.03%  "

$ws.Range("D38").Value = "'0.2692"
$ws.Range("E38").Value = "  -3.17%  "

$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").Value = "'10.36"
$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("E41").Value = "  -7.90%  "

$ws.Range("D42").Value = "'0.7506"
$ws.Range("E42").Value = "  -4.71%  "

$ws.Range("D43").Value = "'13.05"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("D44").Value = "'16.23"
$ws.Range("E44").Value = "  -1.91%  "

$ws.Range("D45").Value = "'0.6894"
$ws.Range("E45").Value = "  -2.82%  "

$ws.Range("D46").Value = "'2.459"
$ws.Range("E46").Value = "  -3.70%  "

$ws.Range("D47").Value = "'4.094"
$ws.Range("E47").Value = "  -2.31%  "

$ws.Range("D48").Value = "'1.002"
$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("D49").Value = "'0.08367"
$ws.Range("E49").Value = "  -3.33%  "

$ws.Range("D50").Value = "'1.264"
$ws.Range("E50").Value = "  -4.83%  "

$ws.Range("D51").Value = "'133.70"
$ws.Range("E51").Value = "  -2.81%  "
